$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "白银有色"
$ws.Range("B2").Value = "白银有色"
$ws.Range("C2").Value = "蓝色光标"

$ws.Range("A3").Value = "利欧股份"
$ws.Range("B3").Value = "湖南白银"
$ws.Range("C3").Value = "白银有色"

$ws.Range("A4").Value = "蓝色光标"
$ws.Range("B4").Value = "利欧股份"
$ws.Range("C4").Value = "巨力索具"

$ws.Range("A5").Value = "湖南白银"
$ws.Range("B5").Value = "铜陵有色"
$ws.Range("C5").Value = "利欧股份"

$ws.Range("A6").Value = "中国黄金"
$ws.Range("B6").Value = "紫金矿业"
$ws.Range("C6").Value = "湖南白银"

$ws.Range("A7").Value = "华天科技"
$ws.Range("B7").Value = "中国铝业"
$ws.Range("C7").Value = "洛阳钼业"

$ws.Range("A8").Value = "铜陵有色"
$ws.Range("B8").Value = "中国黄金"
$ws.Range("C8").Value = "沃尔核材"

$ws.Range("A9").Value = "紫金矿业"
$ws.Range("B9").Value = "华天科技"
$ws.Range("C9").Value = "铜陵有色"

$ws.Range("A10").Value = "巨力索具"
$ws.Range("B10").Value = "蓝色光标"
$ws.Range("C10").Value = "紫金矿业"

$ws.Range("A11").Value = "洛阳钼业"
$ws.Range("B11").Value = "平潭发展"
$ws.Range("C11").Value = "洲际油气"

$ws.Range("A12").Value = "中国铝业"
$ws.Range("B12").Value = "湖南黄金"
$ws.Range("C12").Value = "中国铝业"

$ws.Range("A13").Value = "湖南黄金"
$ws.Range("B13").Value = "洛阳钼业"
$ws.Range("C13").Value = "平潭发展"

$ws.Range("A14").Value = "平潭发展"
$ws.Range("B14").Value = "巨力索具"
$ws.Range("C14").Value = "航天电子"

$ws.Range("A15").Value = "洲际油气"
$ws.Range("B15").Value = "康强电子"
$ws.Range("C15").Value = "浙文互联"

$ws.Range("A16").Value = "航天电子"
$ws.Range("B16").Value = "晓程科技"
$ws.Range("C16").Value = "中国黄金"

$ws.Range("A17").Value = "康强电子"
$ws.Range("B17").Value = "云南铜业"
$ws.Range("C17").Value = "工业富联"

$ws.Range("A18").Value = "晓程科技"
$ws.Range("B18").Value = "北方铜业"
$ws.Range("C18").Value = "华天科技"

$ws.Range("A19").Value = "特  力Ａ"
$ws.Range("B19").Value = "洲际油气"
$ws.Range("C19").Value = "再升科技"

$ws.Range("A20").Value = "黄河旋风"
$ws.Range("B20").Value = "航天电子"
$ws.Range("C20").Value = "康强电子"

$ws.Range("A21").Value = "首都在线"
$ws.Range("B21").Value = "四川黄金"
$ws.Range("C21").Value = "晓程科技"
